$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("C3").Value = "[0, Andre Lucca-Acionamentos-2A, 0,"
$ws.Range("D3").Value = "Josivaldo Ferreira-Circuitos Elétricos 2"

# Row 4
$ws.Range("B4").Value = "-"
$ws.Range("C4").Value = "[0, Andre Lucca-Acionamentos-2A, 0,"
$ws.Range("D4").Value = "Josivaldo Ferreira-Circuitos Elétricos 2"
$ws.Range("F4").Value = "Lucas Ferreira-Sistemas digitais"

# Row 6
$ws.Range("C6").Value = "-"
$ws.Range("D6").Value = "Josivaldo Ferreira-Circuitos Elétricos 2"
$ws.Range("E6").Value = "-"
$ws.Range("F6").Value = "Lucas Ferreira-Sistemas digitais"

# Row 7
$ws.Range("D7").Value = "-"
